# Update the "Haba" price/date records for rows 3-15 (row 2 and row 7 are
# unchanged). Only columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# move; everything else in each row stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    3  = @{ D = 44284; J = 1500; K = 800;  L = 850;  M = 825;  P = 825  }
    4  = @{ D = 44550; J = 1300; K = 1000; L = 1200; M = 1100; P = 1100 }
    5  = @{ D = 44484; J = 900;  K = 750;  L = 800;  M = 775;  P = 775  }
    6  = @{ D = 44243; J = 1200; K = 1200; L = 1300; M = 1250; P = 1250 }
    8  = @{ D = 44453; J = 1000; K = 800;  L = 900;  M = 850;  P = 850  }
    9  = @{ D = 44449; J = 1300; K = 900;  L = 950;  M = 925;  P = 925  }
    10 = @{ D = 44341; J = 1300; K = 900;  L = 1000; M = 950;  P = 950  }
    11 = @{ D = 44442; J = 1250; K = 850;  L = 900;  M = 875;  P = 875  }
    12 = @{ D = 44291; J = 1000; K = 1000; L = 1200; M = 1100; P = 1100 }
    13 = @{ D = 44175; J = 1600; K = 1000; L = 1200; M = 1100; P = 1100 }
    14 = @{ D = 44455; J = 1100; K = 900;  L = 1000; M = 950;  P = 950  }
    15 = @{ D = 44229; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
